# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 418 (pushing existing rows 418-440
# down to 419-441) on the single data sheet, then populate the new row
# with the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 418; this shifts rows
# 418:440 down to 419:441 (formatting/styles carried along by Excel).
$ws.Rows("418:418").Insert()

# Populate the newly inserted row 418 with the new weekly record.
$ws.Cells.Item(418, 1).Value2 = 3
$ws.Cells.Item(418, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(418, 3).Value2 = "Coquimbo"
$ws.Cells.Item(418, 4).Value2 = 45041
$ws.Cells.Item(418, 5).Value2 = 5
$ws.Cells.Item(418, 6).Value2 = 100112001
$ws.Cells.Item(418, 7).Value2 = "Berenjena"
$ws.Cells.Item(418, 8).Value2 = "Sin especificar"
$ws.Cells.Item(418, 9).Value2 = "Primera"
$ws.Cells.Item(418, 10).Value2 = 40
$ws.Cells.Item(418, 11).Value2 = 8000
$ws.Cells.Item(418, 12).Value2 = 8000
$ws.Cells.Item(418, 13).Value2 = 8000
$ws.Cells.Item(418, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(418, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(418, 16).Value2 = 133
$ws.Cells.Item(418, 17).Value2 = 60
$ws.Cells.Item(418, 18).Value2 = "Hortaliza"
